$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume-change (E) columns to match latest scrape.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.644.40"
$ws.Range("E2").Value = "  +6.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.623.67"
$ws.Range("E3").Value = "  +5.78%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.71"
$ws.Range("E5").Value = "  +3.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "191.89"
$ws.Range("E6").Value = "  +6.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.650"
$ws.Range("E7").Value = "  +2.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.617.72"
$ws.Range("E8").Value = "  +5.82%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.180"
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("E11").Value = "  +3.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.35"
$ws.Range("E12").Value = "  +6.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000295"
$ws.Range("E13").Value = "  +5.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.80"
$ws.Range("E14").Value = "  +5.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.200.26"
$ws.Range("E15").Value = "  +5.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.624.23"
$ws.Range("E16").Value = "  +6.24%  "
$ws.Range("E17").Value = "  +5.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.620.35"
$ws.Range("E18").Value = "  +6.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.54"
$ws.Range("E19").Value = "  +4.38%  "
$ws.Range("E21").Value = "  +4.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "498.14"
$ws.Range("E22").Value = "  +7.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.26"
$ws.Range("E23").Value = "  +17.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.38"
$ws.Range("E24").Value = "  +7.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.50"
$ws.Range("E25").Value = "  +8.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.10"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.12"
$ws.Range("E27").Value = "  +5.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.21"
$ws.Range("E28").Value = "  +3.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.45"
$ws.Range("E29").Value = "  +6.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.43"
$ws.Range("E30").Value = "  +3.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.56"
$ws.Range("E31").Value = "  +8.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "623.16"
$ws.Range("E32").Value = "  +6.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.24"
$ws.Range("E33").Value = "  +5.61%  "
$ws.Range("E34").Value = "  +7.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "65.16"
$ws.Range("E35").Value = "  +4.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0830"
$ws.Range("E36").Value = "  +8.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.410"
$ws.Range("E37").Value = "  +6.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.16"
$ws.Range("E38").Value = "  +4.37%  "
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  +1.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.323.81"
$ws.Range("E42").Value = "  +6.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.08"
$ws.Range("E43").Value = "  +4.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0448"
$ws.Range("E44").Value = "  +5.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.69"
$ws.Range("E45").Value = "  +6.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.30"
$ws.Range("E46").Value = "  +2.97%  "
$ws.Range("E47").Value = "  +2.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.15"
$ws.Range("E48").Value = "  +6.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.71"
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.32"
$ws.Range("E50").Value = "  +4.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.21%  "
